# Fixed update to excel issue
# - Rename "Requested quantity" header to "Weekly_PO_Qty" on "Weekly Quantity" sheet
# - Rename "Requested quantity" header to "Monthly_PO_Qty" on "Monthly Trend" sheet
# - Add a new "PO Forecast" sheet with forecast data (ds, PO_Forecast, yhat_lower, yhat_upper)

$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

# Rename the "Requested quantity" headers on the existing sheets.
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" sheet after the last existing sheet, using a copy
# of an existing sheet as a template so sheet-level formatting (outline
# properties, default column width, page margins, header/date cell styles)
# matches the rest of the workbook.
$wsMonthly.Copy($null, $wsMonthly)
$wsForecast = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast.Name = "PO Forecast"

# Header row (A1/B1 inherit the existing bold/centered header style; extend
# that same style to the two new columns before overwriting their text).
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("A1:B1").Copy()
$wsForecast.Range("C1:D1").PasteSpecial(-4122)
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Forecast data rows (dates in column A use the same date-time number format
# as the other sheets' date columns, extended down to row 23).
$data = @(
  @(45333.99999999999, 26, -0.3612331784188433, 53.23875095437787),
  @(45340.99999999999, 25, -3.50469053820449, 53.04373918714737),
  @(45375.99999999999, 21, -7.688249747937048, 49.30620356655725),
  @(45389.99999999999, 20, -8.251577153127812, 48.32978573164073),
  @(45396.99999999999, 19, -9.58524679264651, 46.79673166075153),
  @(45410.99999999999, 17, -12.06774526646888, 44.57958347653568),
  @(45417.99999999999, 17, -11.94106532794404, 42.80295580472231),
  @(45424.99999999999, 16, -9.048226446600328, 43.593642940664),
  @(45431.99999999999, 15, -12.33226326913667, 44.11670616046816),
  @(45438.99999999999, 14, -13.0537483174446, 42.85938650098072),
  @(45445.99999999999, 14, -14.53812827041868, 39.3686175481219),
  @(45452.99999999999, 13, -14.90789126047906, 39.67264439982397),
  @(45459.99999999999, 12, -15.7912448591419, 39.7281610325041),
  @(45487.99999999999, 9, -17.79267399639423, 36.14255456090667),
  @(45494.99999999999, 8, -19.37226877299943, 38.90083609837055),
  @(45501.99999999999, 7, -20.4278986253435, 35.39849496431461),
  @(45508.99999999999, 7, -20.24216489122334, 32.87954814535546),
  @(45515.99999999999, 6, -22.53580634829137, 34.21127193365466),
  @(45522.99999999999, 5, -22.40445733990147, 32.9669417592625),
  @(45529.99999999999, 4, -24.90691585026585, 30.98971211160882),
  @(45536.99999999999, 4, -25.63712628328992, 30.90634825245601),
  @(45543.99999999999, 3, -25.39619627497002, 28.31799552278)
)

$rowCount = $data.Count
$wsForecast.Range("A2").Copy()
$wsForecast.Range("A2:A" + (1 + $rowCount)).PasteSpecial(-4122)

$r = 2
foreach ($row in $data) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r++
}
